$d = $word.ActiveDocument

# Locate the unique "<m>caput mortuum" sequence to scope the edit.
$r = $d.Content
$found = $r.Find.Execute("<m>caput mortuum", $true, $false, $false, $false, $false,
                          $true, 1, $false, "", 0)

if ($found) {
    # Isolate just the "<m>" portion (first 3 characters of the match) and
    # append "<la>" to it -> run becomes "<m><la>".
    $mTagRange = $d.Range($r.Start, $r.Start + 3)
    $mTagRange.InsertAfter("<la>")

    # Re-find "caput mortuum" (now shifted) and append "</la>" to it
    # -> run becomes "caput mortuum</la>".
    $capRange = $d.Content
    $capRange.Find.Execute("caput mortuum", $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)
    $capRange.InsertAfter("</la>")
}
